$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Target state: A1 = "username", A2 = "Xiao(xiao)" only.
# Set A2 to the value that should remain, then delete the now-extra rows 3:5
# (which held "gggggg(gg)", "llll(ll)", and the original "Xiao(xiao)").
$ws.Range("A2").Value = "Xiao(xiao)"
$ws.Range("A3:A5").EntireRow.Delete()
